$wb = $excel.ActiveWorkbook

# Hunk 0 sheet=ALC row=8
$ws = $wb.Worksheets.Item(1)
$ws.Range("H8").Value = 298.42856
$ws.Range("I8").Value = 298.42856
$ws.Range("K8").Value = 895.28568
$ws.Range("M8").Value = -756.28568

# Hunk 1 sheet=ALC row=98
$ws.Range("H98").Value = 157.22223
$ws.Range("I98").Value = 157.22223
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 157.22223
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 1340.77777
$ws.Range("N98").ClearContents() | Out-Null

# Hunk 2 sheet=ALC row=112
$ws.Range("H112").Value = 8137.136
$ws.Range("I112").Value = 260
$ws.Range("J112").Value = 10453.941
$ws.Range("K112").Value = 780
$ws.Range("L112").Value = 31361.823
$ws.Range("M112").Value = 328
$ws.Range("N112").Value = -33577.823

# Hunk 3 sheet=ALC row=122
$ws.Range("H122").Value = 157.22223
$ws.Range("I122").Value = 157.22223
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 471.66669
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 1978.33331
$ws.Range("N122").ClearContents() | Out-Null

# Hunk 4 sheet=ALC row=138
$ws.Range("H138").Value = 1759.6129
$ws.Range("I138").Value = 1371.0588
$ws.Range("J138").Value = 2231.4285
$ws.Range("K138").Value = 4113.1764
$ws.Range("L138").Value = 6694.2855
$ws.Range("M138").Value = 1026.8236
$ws.Range("N138").Value = -16974.2855

# Hunk 5 sheet=ARM row=45
$ws = $wb.Worksheets.Item(2)
$ws.Range("H45").Value = 3438.389
$ws.Range("I45").Value = 3756.2222
$ws.Range("J45").Value = 3120.5557
$ws.Range("K45").Value = 3756.2222
$ws.Range("L45").Value = 3120.5557
$ws.Range("M45").Value = -3379.2222
$ws.Range("N45").Value = -3874.5557

# Hunk 6 sheet=ARM row=97
$ws.Range("H97").Value = 1392.5454
$ws.Range("I97").Value = 1477.5555
$ws.Range("J97").Value = 1010
$ws.Range("K97").Value = 1477.5555
$ws.Range("L97").Value = 1010
$ws.Range("M97").Value = -981.5554999999999
$ws.Range("N97").Value = -2002

# Hunk 7 sheet=ARM row=107
$ws.Range("H107").Value = 36276
$ws.Range("J107").Value = 36276
$ws.Range("L107").Value = 36276
$ws.Range("N107").Value = -43956

# Hunk 8 sheet=ARM row=109
$ws.Range("H109").Value = 44500
$ws.Range("J109").Value = 44500
$ws.Range("L109").Value = 44500
$ws.Range("N109").Value = -47274

# Hunk 9 sheet=ARM row=110
$ws.Range("H110").Value = 2493.8
$ws.Range("I110").Value = 2156.3333
$ws.Range("J110").Value = 3000
$ws.Range("K110").Value = 2156.3333
$ws.Range("L110").Value = 3000
$ws.Range("M110").Value = -111.3332999999998
$ws.Range("N110").Value = -7090

# Hunk 10 sheet=ARM row=112
$ws.Range("H112").Value = 49483.727
$ws.Range("J112").Value = 49483.727
$ws.Range("L112").Value = 49483.727
$ws.Range("N112").Value = -52437.727

# Hunk 11 sheet=ARM row=114
$ws.Range("H114").Value = 38199
$ws.Range("J114").Value = 38199
$ws.Range("L114").Value = 38199
$ws.Range("N114").Value = -46877

# Hunk 12 sheet=ARM row=122
$ws.Range("H122").Value = 168218.5
$ws.Range("I122").Value = 200862.2
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 602586.6000000001
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -600136.6000000001
$ws.Range("N122").Value = -19900

# Hunk 13 sheet=BSM row=107
$ws = $wb.Worksheets.Item(3)
$ws.Range("H107").Value = 38171.297
$ws.Range("I107").Value = 46537.5
$ws.Range("J107").Value = 1360
$ws.Range("K107").Value = 46537.5
$ws.Range("L107").Value = 1360
$ws.Range("M107").Value = -44617.5
$ws.Range("N107").Value = -5200

# Hunk 14 sheet=CRP row=31
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 5594.5
$ws.Range("I31").Value = 1121.7354
$ws.Range("J31").Value = 13198.2
$ws.Range("K31").Value = 1121.7354
$ws.Range("L31").Value = 13198.2
$ws.Range("M31").Value = -826.7354
$ws.Range("N31").Value = -13788.2

# Hunk 15 sheet=CRP row=34
$ws.Range("H34").Value = 5594.5
$ws.Range("I34").Value = 1121.7354
$ws.Range("J34").Value = 13198.2
$ws.Range("K34").Value = 1121.7354
$ws.Range("L34").Value = 13198.2
$ws.Range("M34").Value = -919.7354
$ws.Range("N34").Value = -13602.2

# Hunk 16 sheet=CRP row=107
$ws.Range("H107").Value = 4167540
$ws.Range("I107").Value = 6945211
$ws.Range("J107").Value = 1033.3334
$ws.Range("K107").Value = 6945211
$ws.Range("L107").Value = 1033.3334
$ws.Range("M107").Value = -6943291
$ws.Range("N107").Value = -4873.3334

# Hunk 17 sheet=CRP row=122
$ws.Range("H122").Value = 1649.3695
$ws.Range("I122").Value = 1531.1072
$ws.Range("J122").Value = 1833.3334
$ws.Range("K122").Value = 4593.321599999999
$ws.Range("L122").Value = 5500.0002
$ws.Range("M122").Value = -2143.321599999999
$ws.Range("N122").Value = -10400.0002

# Hunk 18 sheet=CRP row=132
$ws.Range("H132").Value = 5954622
$ws.Range("I132").Value = 1981.8572
$ws.Range("J132").Value = 11907262
$ws.Range("K132").Value = 5945.571599999999
$ws.Range("L132").Value = 35721786
$ws.Range("M132").Value = -3415.571599999999
$ws.Range("N132").Value = -35726846

# Hunk 19 sheet=CUL row=5
$ws = $wb.Worksheets.Item(5)
$ws.Range("H5").Value = 588.6667
$ws.Range("I5").Value = 419.22726
$ws.Range("K5").Value = 1257.68178
$ws.Range("M5").Value = -1145.68178

# Hunk 20 sheet=CUL row=122
$ws.Range("H122").Value = 7965.143
$ws.Range("J122").Value = 15215.429
$ws.Range("L122").Value = 136938.861
$ws.Range("N122").Value = -141838.861

# Hunk 21 sheet=CUL row=123
$ws.Range("H123").Value = 5790.4287
$ws.Range("I123").Value = 2000
$ws.Range("J123").Value = 6422.1665
$ws.Range("K123").Value = 6000
$ws.Range("L123").Value = 19266.4995
$ws.Range("M123").Value = -3550
$ws.Range("N123").Value = -24166.4995

# Hunk 22 sheet=CUL row=131
$ws.Range("H131").Value = 1017.7
$ws.Range("J131").Value = 1127.5883
$ws.Range("L131").Value = 3382.7649
$ws.Range("N131").Value = -13462.7649

# Hunk 23 sheet=CUL row=135
$ws.Range("H135").Value = 588.6667
$ws.Range("I135").Value = 419.22726
$ws.Range("K135").Value = 3773.04534
$ws.Range("M135").Value = -1238.04534

# Hunk 24 sheet=CUL row=136
$ws.Range("H136").Value = 2999.875

# Hunk 25 sheet=CUL row=137
$ws.Range("H137").Value = 13902271
$ws.Range("I137").Value = 27797696
$ws.Range("J137").Value = 6846.6665
$ws.Range("K137").Value = 83393088
$ws.Range("L137").Value = 20539.9995
$ws.Range("M137").Value = -83387988
$ws.Range("N137").Value = -30739.9995

# Hunk 26 sheet=GSM row=122
$ws = $wb.Worksheets.Item(6)
$ws.Range("H122").Value = 1735.625
$ws.Range("I122").Value = 1735.625
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5206.875
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2756.875
$ws.Range("N122").ClearContents() | Out-Null

# Hunk 27 sheet=LTW row=110
$ws = $wb.Worksheets.Item(7)
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents() | Out-Null

# Hunk 28 sheet=LTW row=122
$ws.Range("H122").Value = 2652.1333
$ws.Range("I122").Value = 2394.0588
$ws.Range("K122").Value = 7182.176399999999
$ws.Range("M122").Value = -4732.176399999999

# Hunk 29 sheet=WVR row=122
$ws = $wb.Worksheets.Item(8)
$ws.Range("H122").Value = 1908.4419
$ws.Range("I122").Value = 1470.9656
$ws.Range("J122").Value = 2814.6428
$ws.Range("K122").Value = 4412.8968
$ws.Range("L122").Value = 8443.928400000001
$ws.Range("M122").Value = -1962.8968
$ws.Range("N122").Value = -13343.9284
